$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.400.89'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '2.645.88'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.87'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '2.644.17'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.145'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.34%  '
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.355'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.13'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000193'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.25%  '
$ws.Range('D16').Value = '3.123.54'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').Value = '68.244.70'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '2.652.95'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '364.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.91'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.84'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('E28').Value = '  +2.08%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '573.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.87'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('E35').Value = '  +3.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.38'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  +1.32%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.95%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₆0337'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('E45').Value = '  +3.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.42'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.77'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('E51').Value = '  +1.49%  '
